# Add a "Save" column (H) to the s_vals sheet, matching the existing
# header style used by the other headers (e.g. G1), with zero-valued
# data cells beneath it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the neighboring header cell (G1) onto H1 so the
# new header picks up the same bold/centered/bordered style, then set its
# text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Data rows: plain numeric zeros, no special style (matches F/G columns'
# un-styled data cells).
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0

$excel.CutCopyMode = $false
